$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.421.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.99%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.433.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.434.48"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.86%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +2.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.65%  "

$ws.Range("E11").Value = "  +3.39%  "

$ws.Range("E12").Value = "  +1.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.020.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.123"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.45%  "

$ws.Range("E16").Value = "  +2.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.438.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.531.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.03%  "

$ws.Range("E19").Value = "  +8.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.08%  "

$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "395.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.39%  "

$ws.Range("E23").Value = "  +3.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.49%  "

$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("E26").Value = "  +0.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000123"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.573.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.180"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.91%  "

$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.35%  "

$ws.Range("E33").Value = "  -8.20%  "

$ws.Range("E34").Value = "  +2.52%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.48%  "

$ws.Range("E37").Value = "  +3.71%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.461.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.24%  "

$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "167.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.80%  "

$ws.Range("E42").Value = "  +3.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.799"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.64%  "

$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.593.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.63%  "

$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("E51").Value = "  +2.79%  "
